$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.642.12"
$ws.Range("E2").Value = "  +2.62%  "
$ws.Range("D3").Value = "1.913.85"
$ws.Range("E3").Value = "  +5.63%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'314.05"
$ws.Range("E5").Value = "  +1.52%  "
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.5055"
$ws.Range("E7").Value = "  +2.36%  "
$ws.Range("D8").Value = "'0.3977"
$ws.Range("E8").Value = "  +2.52%  "
$ws.Range("D9").Value = "'0.09680"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").Value = "'1.166"
$ws.Range("E10").Value = "  +5.84%  "
$ws.Range("D11").Value = "'41.75"
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("D12").Value = "'6.554"
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("D14").Value = "1.916.05"
$ws.Range("E14").Value = "  +5.82%  "
$ws.Range("D15").Value = "'7.592"
$ws.Range("E15").Value = "  +4.04%  "
$ws.Range("D16").Value = "'0.9996"
$ws.Range("D17").Value = "'0.00001138"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "'94.21"
$ws.Range("E18").Value = "  +1.81%  "
$ws.Range("D19").Value = "'0.06631"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("E20").Value = "  +5.60%  "
$ws.Range("D21").Value = "'0.9995"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "'6.297"
$ws.Range("E22").Value = "  +6.10%  "
$ws.Range("D23").Value = "28.695.81"
$ws.Range("E23").Value = "  +2.62%  "
$ws.Range("D24").Value = "'11.46"
$ws.Range("E24").Value = "  +2.84%  "
$ws.Range("D25").Value = "'2.283"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("D26").Value = "'2.802"
$ws.Range("E26").Value = "  +17.18%  "
$ws.Range("D27").Value = "2.132.37"
$ws.Range("D28").Value = "'21.45"
$ws.Range("E28").Value = "  +4.07%  "
$ws.Range("D29").Value = "'159.37"
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("D30").Value = "'128.75"
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("D31").Value = "'1.119"
$ws.Range("E31").Value = "  +7.54%  "
$ws.Range("D32").Value = "'0.1077"
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("D33").Value = "'5.743"
$ws.Range("E33").Value = "  +2.99%  "
$ws.Range("D34").Value = "'3.625"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").Value = "'9.815"
$ws.Range("E35").Value = "  +8.87%  "
$ws.Range("D36").Value = "'0.06827"
$ws.Range("E36").Value = "  +1.00%  "
$ws.Range("D37").Value = "'0.02444"
$ws.Range("E37").Value = "  +5.31%  "
$ws.Range("D38").Value = "'0.2222"
$ws.Range("E38").Value = "  +4.24%  "
$ws.Range("D39").Value = "'5.120"
$ws.Range("E39").Value = "  +3.63%  "
$ws.Range("D40").Value = "'11.66"
$ws.Range("E40").Value = "  +3.61%  "
$ws.Range("D41").Value = "'0.6437"
$ws.Range("E41").Value = "  +4.01%  "
$ws.Range("D42").Value = "'1.197"
$ws.Range("E42").Value = "  +4.61%  "
$ws.Range("D43").Value = "'0.9991"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "'13.76"
$ws.Range("E44").Value = "  +5.01%  "
$ws.Range("D45").Value = "'0.6105"
$ws.Range("E45").Value = "  +4.08%  "
$ws.Range("D46").Value = "'1.288"
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("D47").Value = "'3.663"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").Value = "'2.046"
$ws.Range("E48").Value = "  +5.74%  "
$ws.Range("D49").Value = "'124.94"
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("D50").Value = "'1.212"
$ws.Range("E50").Value = "  +3.21%  "
$ws.Range("D51").Value = "'78.50"
$ws.Range("E51").Value = "  +6.78%  "
